$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 178, pushing the existing rows 178:251 down to 179:252.
$ws.Rows.Item(178).Insert()

# Populate the newly inserted row 178 with its data. Columns that are
# identical to the surrounding records (A,B,C,E,F,G,H,I,N,Q,R) are copied
# from the row immediately below (now row 179, the former row 178), while
# the changed columns (D,J,K,L,M,O,P) get their new values.
$ws.Cells.Item(178, 1).Value2 = $ws.Cells.Item(179, 1).Value2    # A: Mercado ID
$ws.Cells.Item(178, 2).Value2 = $ws.Cells.Item(179, 2).Value2    # B: Mercado
$ws.Cells.Item(178, 3).Value2 = $ws.Cells.Item(179, 3).Value2    # C: Region
$ws.Cells.Item(178, 4).Value2 = 44636                            # D: Fecha
$ws.Cells.Item(178, 5).Value2 = $ws.Cells.Item(179, 5).Value2    # E: Codreg
$ws.Cells.Item(178, 6).Value2 = $ws.Cells.Item(179, 6).Value2    # F: Categoria ID
$ws.Cells.Item(178, 7).Value2 = $ws.Cells.Item(179, 7).Value2    # G: Categoria
$ws.Cells.Item(178, 8).Value2 = $ws.Cells.Item(179, 8).Value2    # H: Variedad
$ws.Cells.Item(178, 9).Value2 = $ws.Cells.Item(179, 9).Value2    # I: Calidad
$ws.Cells.Item(178, 10).Value2 = 30                              # J: Volumen
$ws.Cells.Item(178, 11).Value2 = 12000                           # K: Precio minimo
$ws.Cells.Item(178, 12).Value2 = 12000                           # L: Precio maximo
$ws.Cells.Item(178, 13).Value2 = 12000                           # M: Precio promedio ponderado
$ws.Cells.Item(178, 14).Value2 = $ws.Cells.Item(179, 14).Value2  # N: Unidad de comercializacion
$ws.Cells.Item(178, 15).Value2 = "Región del Maule"              # O: Origen
$ws.Cells.Item(178, 16).Value2 = 200                             # P: Precio $/Kg
$ws.Cells.Item(178, 17).Value2 = $ws.Cells.Item(179, 17).Value2  # Q: Kg o Unidades
$ws.Cells.Item(178, 18).Value2 = $ws.Cells.Item(179, 18).Value2  # R: Clasificacion
